$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 3).Value = 45
$ws.Cells.Item(1, 5).Value = 360
$ws.Cells.Item(2, 3).Value = 45
$ws.Cells.Item(2, 5).Value = 360
$ws.Cells.Item(3, 3).Value = 45
$ws.Cells.Item(3, 5).Value = 315
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 5).Value = 120
$ws.Cells.Item(5, 3).Value = 20
$ws.Cells.Item(5, 5).Value = 120
$ws.Cells.Item(6, 3).Value = 20
$ws.Cells.Item(6, 5).Value = 100
$ws.Cells.Item(7, 3).Value = 5
$ws.Cells.Item(7, 5).Value = 50
$ws.Cells.Item(8, 3).Value = 10
$ws.Cells.Item(8, 5).Value = 50
$ws.Cells.Item(9, 3).Value = 5
$ws.Cells.Item(9, 5).Value = 35
$ws.Cells.Item(10, 3).Value = 60
$ws.Cells.Item(10, 5).Value = 360
$ws.Cells.Item(11, 3).Value = 50
$ws.Cells.Item(11, 5).Value = 350
$ws.Cells.Item(12, 3).Value = 50
$ws.Cells.Item(12, 5).Value = 300
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(13, 5).Value = 160
$ws.Cells.Item(14, 3).Value = 20
$ws.Cells.Item(14, 5).Value = 140
$ws.Cells.Item(15, 3).Value = 20
$ws.Cells.Item(15, 5).Value = 120
$ws.Cells.Item(16, 3).Value = 70
$ws.Cells.Item(16, 5).Value = 700
$ws.Cells.Item(17, 3).Value = 70
$ws.Cells.Item(17, 5).Value = 700
$ws.Cells.Item(18, 3).Value = 70
$ws.Cells.Item(18, 5).Value = 420
$ws.Cells.Item(19, 3).Value = 45
$ws.Cells.Item(19, 5).Value = 450
$ws.Cells.Item(20, 3).Value = 45
$ws.Cells.Item(20, 5).Value = 450
$ws.Cells.Item(21, 3).Value = 45
$ws.Cells.Item(21, 5).Value = 360
$ws.Cells.Item(22, 3).Value = 35
$ws.Cells.Item(22, 5).Value = 280
$ws.Cells.Item(23, 3).Value = 35
$ws.Cells.Item(23, 5).Value = 210
$ws.Cells.Item(24, 3).Value = 20
$ws.Cells.Item(24, 5).Value = 160
$ws.Cells.Item(25, 3).Value = 35
$ws.Cells.Item(25, 5).Value = 210
$ws.Cells.Item(26, 3).Value = 30
$ws.Cells.Item(26, 5).Value = 180
$ws.Cells.Item(27, 3).Value = 25
$ws.Cells.Item(27, 5).Value = 150
$ws.Cells.Item(28, 3).Value = 25
$ws.Cells.Item(28, 5).Value = 150
$ws.Cells.Item(29, 3).Value = 20
$ws.Cells.Item(29, 5).Value = 120
$ws.Cells.Item(30, 3).Value = 20
$ws.Cells.Item(30, 5).Value = 120
$ws.Cells.Item(31, 3).Value = 20
$ws.Cells.Item(31, 5).Value = 120
$ws.Cells.Item(32, 3).Value = 12.5
$ws.Cells.Item(32, 5).Value = 75
$ws.Cells.Item(33, 3).Value = 12.5
$ws.Cells.Item(33, 5).Value = 75
$ws.Cells.Item(34, 3).Value = 12.5
$ws.Cells.Item(34, 5).Value = 50
$ws.Cells.Item(35, 3).Value = 50
$ws.Cells.Item(35, 5).Value = 400
$ws.Cells.Item(36, 3).Value = 50
$ws.Cells.Item(36, 5).Value = 400
$ws.Cells.Item(37, 3).Value = 50
$ws.Cells.Item(37, 5).Value = 400
$ws.Cells.Item(38, 3).Value = 15
$ws.Cells.Item(38, 5).Value = 105
$ws.Cells.Item(39, 3).Value = 15
$ws.Cells.Item(39, 5).Value = 105
$ws.Cells.Item(40, 3).Value = 15
$ws.Cells.Item(40, 5).Value = 105
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(44, 3).Value = 25
$ws.Cells.Item(44, 5).Value = 250
$ws.Cells.Item(45, 3).Value = 35
$ws.Cells.Item(45, 5).Value = 350
$ws.Cells.Item(46, 3).Value = 45
$ws.Cells.Item(46, 5).Value = 450
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(49, 2).Value = "Single Leg Standing Calf Raise (Barbell)"
$ws.Cells.Item(49, 4).Value = 10
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(50, 3).Value = 35
$ws.Cells.Item(50, 5).Value = 210
$ws.Cells.Item(51, 3).Value = 35
$ws.Cells.Item(51, 5).Value = 210
$ws.Cells.Item(52, 3).Value = 35
$ws.Cells.Item(52, 5).Value = 245
$ws.Cells.Item(53, 3).Value = 30
$ws.Cells.Item(53, 5).Value = 180
$ws.Cells.Item(54, 3).Value = 30
$ws.Cells.Item(54, 5).Value = 180
$ws.Cells.Item(55, 3).Value = 20
$ws.Cells.Item(55, 5).Value = 120
$ws.Cells.Item(56, 3).Value = 60
$ws.Cells.Item(56, 5).Value = 300
$ws.Cells.Item(57, 3).Value = 50
$ws.Cells.Item(57, 5).Value = 250
$ws.Cells.Item(58, 3).Value = 40
$ws.Cells.Item(58, 5).Value = 240
$ws.Cells.Item(59, 3).Value = 20
$ws.Cells.Item(59, 5).Value = 200
$ws.Cells.Item(60, 3).Value = 20
$ws.Cells.Item(60, 5).Value = 200
$ws.Cells.Item(61, 3).Value = 20
$ws.Cells.Item(61, 5).Value = 160
$ws.Cells.Item(62, 3).Value = 25
$ws.Cells.Item(62, 5).Value = 250
$ws.Cells.Item(63, 3).Value = 35
$ws.Cells.Item(63, 5).Value = 350
$ws.Cells.Item(64, 3).Value = 35
$ws.Cells.Item(64, 5).Value = 385
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 3).Value = 20
$ws.Cells.Item(66, 5).Value = 200
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 5).Value = 0
$ws.Cells.Item(68, 3).Value = 45
$ws.Cells.Item(68, 5).Value = 360
$ws.Cells.Item(69, 3).Value = 45
$ws.Cells.Item(69, 5).Value = 360
$ws.Cells.Item(70, 3).Value = 45
$ws.Cells.Item(70, 5).Value = 360
$ws.Cells.Item(71, 3).Value = 80
$ws.Cells.Item(71, 5).Value = 480
$ws.Cells.Item(72, 3).Value = 60
$ws.Cells.Item(72, 5).Value = 480
$ws.Cells.Item(73, 3).Value = 45
$ws.Cells.Item(73, 5).Value = 270
